function Set-CellText {
    param($ws, $ref, $val)
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $val
    $ws.Range($ref).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws 'D2' '58.270.38'
Set-CellText $ws 'E2' '  -1.38%  '
Set-CellText $ws 'D3' '2.584.68'
Set-CellText $ws 'E3' '  -2.06%  '
Set-CellText $ws 'E4' '  +0.07%  '
Set-CellText $ws 'D5' '518.63'
Set-CellText $ws 'E5' '  -1.70%  '
Set-CellText $ws 'D6' '140.84'
Set-CellText $ws 'E6' '  -2.58%  '
Set-CellText $ws 'E7' '  -0.11%  '
Set-CellText $ws 'D8' '0.563'
Set-CellText $ws 'E8' '  -0.97%  '
Set-CellText $ws 'D9' '2.603.34'
Set-CellText $ws 'E9' '  -1.77%  '
Set-CellText $ws 'D10' '6.46'
Set-CellText $ws 'E10' '  -2.72%  '
Set-CellText $ws 'D11' '0.100'
Set-CellText $ws 'E11' '  -3.40%  '
Set-CellText $ws 'D12' '0.330'
Set-CellText $ws 'E12' '  -1.80%  '
Set-CellText $ws 'E13' '  +0.45%  '
Set-CellText $ws 'D14' '3.042.16'
Set-CellText $ws 'E14' '  -2.01%  '
Set-CellText $ws 'D15' '58.241.48'
Set-CellText $ws 'E15' '  -1.34%  '
Set-CellText $ws 'D16' '20.24'
Set-CellText $ws 'E16' '  -3.58%  '
Set-CellText $ws 'D17' '2.593.92'
Set-CellText $ws 'E17' '  -1.01%  '
Set-CellText $ws 'E18' '  -2.82%  '
Set-CellText $ws 'D19' '336.87'
Set-CellText $ws 'E19' '  -1.37%  '
Set-CellText $ws 'D20' '4.30'
Set-CellText $ws 'E20' '  -3.22%  '
Set-CellText $ws 'D21' '10.16'
Set-CellText $ws 'E21' '  -3.60%  '
Set-CellText $ws 'D22' '6.38'
Set-CellText $ws 'E22' '  +0.76%  '
Set-CellText $ws 'D23' '0.997'
Set-CellText $ws 'E23' '  -0.14%  '
Set-CellText $ws 'D24' '65.71'
Set-CellText $ws 'E24' '  +0.87%  '
Set-CellText $ws 'D25' '0.168'
Set-CellText $ws 'E25' '  +0.08%  '
Set-CellText $ws 'D26' '0.401'
Set-CellText $ws 'E26' '  -3.97%  '
Set-CellText $ws 'D27' '0.997'
Set-CellText $ws 'E27' '  -0.09%  '
Set-CellText $ws 'D28' '2.705.20'
Set-CellText $ws 'E28' '  -1.97%  '
Set-CellText $ws 'D29' '6.99'
Set-CellText $ws 'E29' '  -3.45%  '
Set-CellText $ws 'E30' '  -0.06%  '
Set-CellText $ws 'D31' '0.0₃0733'
Set-CellText $ws 'E31' '  -8.11%  '
Set-CellText $ws 'D32' '6.05'
Set-CellText $ws 'E32' '  -6.29%  '
Set-CellText $ws 'E33' '  -2.80%  '
Set-CellText $ws 'D34' '18.71'
Set-CellText $ws 'E34' '  -1.10%  '
Set-CellText $ws 'D35' '149.11'
Set-CellText $ws 'E35' '  -0.68%  '
Set-CellText $ws 'D36' '3.95'
Set-CellText $ws 'E36' '  -5.95%  '
Set-CellText $ws 'D37' '1.12'
Set-CellText $ws 'E37' '  -6.69%  '
Set-CellText $ws 'D38' '0.852'
Set-CellText $ws 'E38' '  -2.31%  '
Set-CellText $ws 'D39' '36.19'
Set-CellText $ws 'E39' '  -0.90%  '
Set-CellText $ws 'E40' '  -1.82%  '
Set-CellText $ws 'D41' '0.826'
Set-CellText $ws 'E41' '  -10.75%  '
Set-CellText $ws 'D42' '3.50'
Set-CellText $ws 'E42' '  -4.31%  '
Set-CellText $ws 'D43' '0.996'
Set-CellText $ws 'E43' '  -0.06%  '
Set-CellText $ws 'D44' '273.17'
Set-CellText $ws 'E44' '  +0.58%  '
Set-CellText $ws 'D45' '0.601'
Set-CellText $ws 'E45' '  -0.29%  '
Set-CellText $ws 'D46' '10.70'
Set-CellText $ws 'E46' '  +0.48%  '
Set-CellText $ws 'D47' '0.0947'
Set-CellText $ws 'E47' '  -2.84%  '
Set-CellText $ws 'D48' '0.0518'
Set-CellText $ws 'E48' '  -3.69%  '
Set-CellText $ws 'D49' '18.55'
Set-CellText $ws 'E49' '  -4.24%  '
Set-CellText $ws 'D50' '1.968.76'
Set-CellText $ws 'E50' '  -3.90%  '
Set-CellText $ws 'D51' '4.55'
Set-CellText $ws 'E51' '  -2.92%  '
